$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-07T07:02:27.221228+00:00"
$ws.Range("K3").Value = "2025-11-07T07:02:27.221264+00:00"
$ws.Range("K4").Value = "2025-11-07T07:02:27.221283+00:00"
$ws.Range("K5").Value = "2025-11-07T07:02:29.962744+00:00"
$ws.Range("K6").Value = "2025-11-07T07:02:29.962772+00:00"
$ws.Range("K7").Value = "2025-11-07T07:02:29.962792+00:00"
$ws.Range("K8").Value = "2025-11-07T07:02:32.691794+00:00"
$ws.Range("K9").Value = "2025-11-07T07:02:35.501998+00:00"
$ws.Range("K10").Value = "2025-11-07T07:02:35.502028+00:00"
$ws.Range("K11").Value = "2025-11-07T07:02:35.502046+00:00"
$ws.Range("K12").Value = "2025-11-07T07:02:38.366426+00:00"
$ws.Range("K13").Value = "2025-11-07T07:02:38.366459+00:00"
$ws.Range("K14").Value = "2025-11-07T07:02:38.366478+00:00"
$ws.Range("K15").Value = "2025-11-07T07:02:38.366503+00:00"
$ws.Range("K16").Value = "2025-11-07T07:02:43.739940+00:00"
$ws.Range("K17").Value = "2025-11-07T07:02:46.487628+00:00"
$ws.Range("K18").Value = "2025-11-07T07:02:48.884537+00:00"
$ws.Range("K19").Value = "2025-11-07T07:02:48.884571+00:00"
$ws.Range("K20").Value = "2025-11-07T07:02:51.716394+00:00"
$ws.Range("K21").Value = "2025-11-07T07:02:54.077386+00:00"
$ws.Range("K22").Value = "2025-11-07T07:02:54.077417+00:00"
$ws.Range("K23").Value = "2025-11-07T07:02:54.077435+00:00"
$ws.Range("K24").Value = "2025-11-07T07:02:56.385463+00:00"
$ws.Range("K25").Value = "2025-11-07T07:02:56.385493+00:00"
$ws.Range("K26").Value = "2025-11-07T07:02:56.385513+00:00"
$ws.Range("K27").Value = "2025-11-07T07:02:59.121532+00:00"
$ws.Range("K28").Value = "2025-11-07T07:02:59.121560+00:00"
$ws.Range("K29").Value = "2025-11-07T07:02:59.121578+00:00"
$ws.Range("K30").Value = "2025-11-07T07:02:59.121593+00:00"
$ws.Range("K31").Value = "2025-11-07T07:02:59.121608+00:00"
$ws.Range("K32").Value = "2025-11-07T07:03:01.467494+00:00"
$ws.Range("K33").Value = "2025-11-07T07:03:04.370984+00:00"
$ws.Range("K34").Value = "2025-11-07T07:03:04.371014+00:00"
$ws.Range("K35").Value = "2025-11-07T07:03:04.371033+00:00"
$ws.Range("K36").Value = "2025-11-07T07:03:07.127198+00:00"
$ws.Range("K37").Value = "2025-11-07T07:03:07.127226+00:00"
$ws.Range("K38").Value = "2025-11-07T07:03:07.127244+00:00"
$ws.Range("K39").Value = "2025-11-07T07:03:09.411035+00:00"
$ws.Range("K40").Value = "2025-11-07T07:03:09.411065+00:00"
$ws.Range("K41").Value = "2025-11-07T07:03:09.411083+00:00"
$ws.Range("K42").Value = "2025-11-07T07:03:09.411098+00:00"
$ws.Range("K43").Value = "2025-11-07T07:03:09.411114+00:00"
$ws.Range("K44").Value = "2025-11-07T07:03:09.411129+00:00"
$ws.Range("K45").Value = "2025-11-07T07:03:09.411144+00:00"
$ws.Range("K46").Value = "2025-11-07T07:03:09.411158+00:00"
$ws.Range("K47").Value = "2025-11-07T07:03:12.361501+00:00"
$ws.Range("K48").Value = "2025-11-07T07:03:12.361531+00:00"
$ws.Range("K49").Value = "2025-11-07T07:03:17.013705+00:00"
$ws.Range("K50").Value = "2025-11-07T07:03:17.013733+00:00"
$ws.Range("K51").Value = "2025-11-07T07:03:19.284191+00:00"
$ws.Range("K52").Value = "2025-11-07T07:03:19.284219+00:00"
